$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-07T21:00:10+00:00"

# --- Sheet "Mapping Table 0": make the three entryRelationship targets more specific ---
$t0 = $wb.Worksheets.Item("Mapping Table 0")
$t0.Range("D16").Value = "FRCDAActe.entryRelationship:frReferenceInterneCirconstances"
$t0.Range("D17").Value = "FRCDAActe.entryRelationship:frReferenceInterneMotifActe"
$t0.Range("D18").Value = "FRCDAActe.entryRelationship:frReferenceInterneDM"

# --- Sheet "Mapping Table 1": fix two source labels, and add a new mapping row ---
$t1 = $wb.Worksheets.Item("Mapping Table 1")
$t1.Range("A10").Value = "FRCDAActe.entryRelationship:frSimpleObservationScores"
$t1.Range("A15").Value = "FRCDAActe.entryRelationship:frReferenceInterneCirconstances"

# Insert a new row 16 (pushing author/priority/bodySite rows down by one)
$t1.Rows.Item(16).Insert()

# Clone formatting from the row above so the new row matches existing styling
$t1.Range("A15:E15").Copy()
$t1.Range("A16:E16").PasteSpecial(-4122)

# Populate the new mapping row
$t1.Range("A16").Value = "FRCDAActe.entryRelationship:frSimpleObservationDifficulte"
$t1.Range("C16").Value = "equivalent"
$t1.Range("D16").Value = "FRProcedureActDocument.extension:difficulte"
